# Weekly update: insert a new price-report row for "Camote" variety "Paine"
# at row 143 (Macroferia Regional de Talca - Zapallo), pushing the
# pre-existing rows 143-147 down to 144-148.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 143; existing rows 143:147 shift down to 144:148.
$ws.Rows.Item(143).Insert()

# Populate the new row 143 with this week's record.
$ws.Range("A143").Value = 5
$ws.Range("B143").Value = "Macroferia Regional de Talca"
$ws.Range("C143").Value = "Maule"
$ws.Range("D143").Value = 44509
$ws.Range("E143").Value = 7
$ws.Range("F143").Value = 100112045
$ws.Range("G143").Value = "Zapallo"
$ws.Range("H143").Value = "Paine"
$ws.Range("I143").Value = "1a (guarda)"
$ws.Range("J143").Value = 2000
$ws.Range("K143").Value = 80
$ws.Range("L143").Value = 80
$ws.Range("M143").Value = 80
$ws.Range("N143").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O143").Value = "Región del Maule"
$ws.Range("P143").Value = 80
$ws.Range("Q143").Value = 1
$ws.Range("R143").Value = "Hortaliza"
